# Agreement migration update:
#  - status column (L) values 81 -> 100 (imported executed / expired agreement)
#  - row 27 & row 77 grow taller (wrapped text now needs more room)
#  - selection / scroll position moves to the bottom of the sheet (last edited row)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- status column: change every "81" in L2:L171 to "100" ---
for ($r = 2; $r -le 171; $r++) {
    $cell = $ws.Cells.Item($r, 12)
    if ($cell.Value2 -eq 81) {
        $cell.Value = 100
    }
}

# --- row height adjustments ---
$ws.Rows.Item(27).RowHeight = 86.25
$ws.Rows.Item(77).RowHeight = 345

# --- move the active selection / viewport towards the end of the sheet ---
$ws.Activate()
$excel.Goto($ws.Range("J164"), $true)
$ws.Range("L171").Select()
